$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4732.0415
$ws.Range("I132").Value = 2189.9048
$ws.Range("K132").Value = 6569.714399999999
$ws.Range("M132").Value = -4039.714399999999
$ws.Range("H133").Value = 109999
$ws.Range("J133").Value = 109999
$ws.Range("L133").Value = 109999
$ws.Range("N133").Value = -120119
$ws.Range("H137").Value = 1534571.4
$ws.Range("I137").Value = 2300402.8
$ws.Range("K137").Value = 6901208.399999999
$ws.Range("M137").Value = -6898658.399999999
$ws.Range("H138").Value = 1768.4546
$ws.Range("J138").Value = 2018.0676
$ws.Range("L138").Value = 6054.2028
$ws.Range("N138").Value = -16334.2028

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7796.367
$ws.Range("I32").Value = 6437.8887
$ws.Range("J32").Value = 13145.375
$ws.Range("K32").Value = 6437.8887
$ws.Range("L32").Value = 13145.375
$ws.Range("M32").Value = -6150.8887
$ws.Range("N32").Value = -13719.375
$ws.Range("H45").Value = 1020.4545
$ws.Range("I45").Value = 1041.5333
$ws.Range("J45").Value = 975.2857
$ws.Range("K45").Value = 1041.5333
$ws.Range("L45").Value = 975.2857
$ws.Range("M45").Value = -664.5333000000001
$ws.Range("N45").Value = -1729.2857
$ws.Range("H61").Value = 6504.3667
$ws.Range("I61").Value = 2913
$ws.Range("J61").Value = 11891.417
$ws.Range("K61").Value = 2913
$ws.Range("L61").Value = 11891.417
$ws.Range("M61").Value = -2701
$ws.Range("N61").Value = -12315.417
$ws.Range("H74").Value = 44376.207
$ws.Range("I74").Value = 48850.723
$ws.Range("J74").Value = 5895.4
$ws.Range("K74").Value = 48850.723
$ws.Range("L74").Value = 5895.4
$ws.Range("M74").Value = -47976.723
$ws.Range("N74").Value = -7643.4
$ws.Range("H77").Value = 44376.207
$ws.Range("I77").Value = 48850.723
$ws.Range("J77").Value = 5895.4
$ws.Range("K77").Value = 244253.615
$ws.Range("L77").Value = 29477
$ws.Range("M77").Value = -239885.615
$ws.Range("N77").Value = -38213
$ws.Range("H102").Value = 1811.05
$ws.Range("I102").Value = 1479.8125
$ws.Range("J102").Value = 3136
$ws.Range("K102").Value = 1479.8125
$ws.Range("L102").Value = 3136
$ws.Range("M102").Value = 142.1875
$ws.Range("N102").Value = -6380
$ws.Range("H122").Value = 5212.9707
$ws.Range("I122").Value = 2984.423
$ws.Range("J122").Value = 12455.75
$ws.Range("K122").Value = 8953.269
$ws.Range("L122").Value = 37367.25
$ws.Range("M122").Value = -6503.269
$ws.Range("N122").Value = -42267.25
$ws.Range("H136").Value = 6504.3667
$ws.Range("I136").Value = 2913
$ws.Range("J136").Value = 11891.417
$ws.Range("K136").Value = 8739
$ws.Range("L136").Value = 35674.251
$ws.Range("M136").Value = -6189
$ws.Range("N136").Value = -40774.251

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1172.2059
$ws.Range("I105").Value = 1172.742
$ws.Range("J105").Value = 1166.6666
$ws.Range("K105").Value = 1172.742
$ws.Range("L105").Value = 1166.6666
$ws.Range("M105").Value = 574.258
$ws.Range("N105").Value = -4660.6666
$ws.Range("H132").Value = 109999
$ws.Range("J132").Value = 109999
$ws.Range("L132").Value = 109999
$ws.Range("N132").Value = -120119

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 48132810
$ws.Range("I31").Value = 626942.1
$ws.Range("K31").Value = 626942.1
$ws.Range("M31").Value = -626647.1
$ws.Range("H34").Value = 48132810
$ws.Range("I34").Value = 626942.1
$ws.Range("K34").Value = 626942.1
$ws.Range("M34").Value = -626740.1
$ws.Range("H58").Value = 2581.08
$ws.Range("I58").Value = 2488.625
$ws.Range("K58").Value = 2488.625
$ws.Range("M58").Value = -2285.625
$ws.Range("H132").Value = 2994.8845
$ws.Range("I132").Value = 1135.5
$ws.Range("K132").Value = 3406.5
$ws.Range("M132").Value = -876.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 4664.3125
$ws.Range("I134").Value = 4864.8667
$ws.Range("J134").Value = 1656
$ws.Range("K134").Value = 14594.6001
$ws.Range("L134").Value = 4968
$ws.Range("M134").Value = -12059.6001
$ws.Range("N134").Value = -10038
$ws.Range("H136").Value = 2581.08
$ws.Range("I136").Value = 2488.625
$ws.Range("K136").Value = 7465.875
$ws.Range("M136").Value = -4915.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 419.66666
$ws.Range("I97").Value = 379.75
$ws.Range("J97").Value = 499.5
$ws.Range("K97").Value = 1139.25
$ws.Range("L97").Value = 1498.5
$ws.Range("M97").Value = -643.25
$ws.Range("N97").Value = -2490.5
$ws.Range("H107").Value = 420.10715
$ws.Range("J107").Value = 424.6
$ws.Range("L107").Value = 1273.8
$ws.Range("N107").Value = -5113.8
$ws.Range("H132").Value = 6131.5835
$ws.Range("I132").Value = 7239.6313
$ws.Range("J132").Value = 1921
$ws.Range("K132").Value = 65156.6817
$ws.Range("L132").Value = 17289
$ws.Range("M132").Value = -62626.6817
$ws.Range("N132").Value = -22349
$ws.Range("H139").Value = 2096.5
$ws.Range("I139").Value = 2096.5
$ws.Range("K139").Value = 6289.5
$ws.Range("M139").Value = -1149.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3983.3
$ws.Range("I97").Value = 3380.25
$ws.Range("J97").Value = 5189.4
$ws.Range("K97").Value = 3380.25
$ws.Range("L97").Value = 5189.4
$ws.Range("M97").Value = -2884.25
$ws.Range("N97").Value = -6181.4
$ws.Range("H132").Value = 24674.803
$ws.Range("I132").Value = 28770.607
$ws.Range("J132").Value = 3786.2
$ws.Range("K132").Value = 86311.821
$ws.Range("L132").Value = 11358.6
$ws.Range("M132").Value = -83781.821
$ws.Range("N132").Value = -16418.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H136").Value = 3482
$ws.Range("I136").Value = 1844.1471
$ws.Range("K136").Value = 5532.4413
$ws.Range("M136").Value = -2982.4413

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1261.7273
$ws.Range("I122").Value = 1087.9
$ws.Range("K122").Value = 3263.7
$ws.Range("M122").Value = -813.7000000000003
$ws.Range("H126").Value = 1780.6154
$ws.Range("I126").Value = 1845.6666
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5536.9998
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -3066.9998
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 19610332
$ws.Range("I132").Value = 25002260
$ws.Range("J132").Value = 3325.7273
$ws.Range("K132").Value = 75006780
$ws.Range("L132").Value = 9977.1819
$ws.Range("M132").Value = -75004250
$ws.Range("N132").Value = -15037.1819

Write-Host "Applied all Lich_Profits updates"